$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values across rows 2-25 (missing-data pattern changes) ---
# --- and fully rewrite rows 26-33 to reflect the re-shuffled/renamed station rows ---
$ws.Cells.Item(3,6).Value = ""
$ws.Cells.Item(4,5).Value = ""
$ws.Cells.Item(5,4).Value = ""
$ws.Cells.Item(6,3).Value = 15.1
$ws.Cells.Item(7,5).Value = ""
$ws.Cells.Item(7,6).Value = 17.24
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(11,4).Value = -15.5
$ws.Cells.Item(12,6).Value = ""
$ws.Cells.Item(13,6).Value = 17.1
$ws.Cells.Item(15,6).Value = ""
$ws.Cells.Item(16,6).Value = 17.34
$ws.Cells.Item(17,5).Value = -7.3
$ws.Cells.Item(17,6).Value = 17.78
$ws.Cells.Item(19,3).Value = 13.2
$ws.Cells.Item(19,4).Value = ""
$ws.Cells.Item(20,6).Value = 17.73
$ws.Cells.Item(21,3).Value = ""
$ws.Cells.Item(22,6).Value = 16.81
$ws.Cells.Item(23,3).Value = 12.2
$ws.Cells.Item(23,4).Value = -13.9
$ws.Cells.Item(24,5).Value = -8.1
$ws.Cells.Item(24,6).Value = ""
$ws.Cells.Item(25,4).Value = -15.5
$ws.Cells.Item(26,1).Value = 'SC 5'
$ws.Cells.Item(26,3).Value = 10.8
$ws.Cells.Item(26,4).Value = -13.8
$ws.Cells.Item(26,5).Value = -5
$ws.Cells.Item(26,6).Value = 17.38
$ws.Cells.Item(27,1).Value = 'SC 101'
$ws.Cells.Item(27,2).Value = -20.4
$ws.Cells.Item(27,3).Value = ""
$ws.Cells.Item(27,4).Value = ""
$ws.Cells.Item(27,5).Value = -10
$ws.Cells.Item(27,6).Value = 17
$ws.Cells.Item(28,1).Value = 'SC 105'
$ws.Cells.Item(28,2).Value = -19.6
$ws.Cells.Item(28,3).Value = 11.1
$ws.Cells.Item(28,4).Value = -13.7
$ws.Cells.Item(28,5).Value = -5.9
$ws.Cells.Item(28,6).Value = ""
$ws.Cells.Item(29,1).Value = 'SC 119'
$ws.Cells.Item(29,3).Value = 11.2
$ws.Cells.Item(29,4).Value = ""
$ws.Cells.Item(29,6).Value = ""
$ws.Cells.Item(30,1).Value = 'SC 120'
$ws.Cells.Item(30,2).Value = -19.7
$ws.Cells.Item(30,3).Value = 11.4
$ws.Cells.Item(30,4).Value = -13.6
$ws.Cells.Item(30,6).Value = ""
$ws.Cells.Item(31,1).Value = 'SC 132'
$ws.Cells.Item(31,2).Value = -18.8
$ws.Cells.Item(31,3).Value = 15.3
$ws.Cells.Item(31,4).Value = -13.7
$ws.Cells.Item(31,5).Value = -8.1
$ws.Cells.Item(31,6).Value = 17.18
$ws.Cells.Item(32,1).Value = 'SC 193'
$ws.Cells.Item(32,2).Value = -19.9
$ws.Cells.Item(32,3).Value = 10.5
$ws.Cells.Item(32,4).Value = -14.7
$ws.Cells.Item(32,5).Value = ""
$ws.Cells.Item(32,6).Value = 17.39
$ws.Cells.Item(33,1).Value = 'SC 232'
$ws.Cells.Item(33,2).Value = -19.5
$ws.Cells.Item(33,3).Value = ""
$ws.Cells.Item(33,4).Value = -14.1
$ws.Cells.Item(33,5).Value = -10.7
$ws.Cells.Item(33,6).Value = ""

# --- Remove the two trailing rows (old SC 193 / SC 232 duplicates), shifting nothing below ---
# This also updates the sheet dimension from A1:F35 to A1:F33
$ws.Rows("34:35").Delete()
